# Weekly data refresh: a new daily/weekly price-report row for
# "Feria Lagunitas de Puerto Montt - Mango" is inserted at the top of the
# data table (row 363), pushing all the existing report rows down by one
# row (old row 363 becomes 364, ..., old row 410 becomes 411).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 363, shifting rows 363:410 down to 364:411.
$ws.Rows("363:363").Insert(-4121)   # -4121 = xlShiftDown

# Populate the new row 363 with the new week's report for this product.
$ws.Range("A363").Value2 = 4
$ws.Range("B363").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C363").Value2 = "Los Lagos"
$ws.Range("D363").Value2 = 45124
$ws.Range("E363").Value2 = 10
$ws.Range("F363").Value2 = "Fruta"
$ws.Range("G363").Value2 = 100108
$ws.Range("H363").Value2 = "Tropicales y subtropicales"
$ws.Range("I363").Value2 = 100108002
$ws.Range("J363").Value2 = "Mango"
$ws.Range("K363").Value2 = "Sin especificar"
$ws.Range("L363").Value2 = "Primera"
$ws.Range("M363").Value2 = 80
$ws.Range("N363").Value2 = 8500
$ws.Range("O363").Value2 = 9000
$ws.Range("P363").Value2 = 8750
$ws.Range("Q363").Value2 = "$/bandeja 4 kilos"
$ws.Range("R363").Value2 = "Brasil"
$ws.Range("S363").Value2 = 2188
$ws.Range("T363").Value2 = 4
